$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns retain their original Text format so
# numeric-looking strings (e.g. "238.75") are not reinterpreted as numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '90.605.83'
$ws.Range('E2').Value = '  -0.74%  '
$ws.Range('D3').Value = '3.153.33'
$ws.Range('E3').Value = '  +1.34%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '238.75'
$ws.Range('E5').Value = '  +9.04%  '
$ws.Range('D6').Value = '643.07'
$ws.Range('E6').Value = '  +3.38%  '
$ws.Range('E7').Value = '  +10.29%  '
$ws.Range('D8').Value = '0.363'
$ws.Range('E8').Value = '  -4.22%  '
$ws.Range('E9').Value = '  +0.03%  '
$ws.Range('D10').Value = '3.133.91'
$ws.Range('E10').Value = '  +0.85%  '
$ws.Range('D11').Value = '0.721'
$ws.Range('E11').Value = '  -0.25%  '
$ws.Range('E12').Value = '  +2.69%  '
$ws.Range('D13').Value = '36.41'
$ws.Range('E13').Value = '  +5.48%  '
$ws.Range('D14').Value = '0.0000248'
$ws.Range('E14').Value = '  -2.65%  '
$ws.Range('D15').Value = '5.64'
$ws.Range('E15').Value = '  +4.50%  '
$ws.Range('D16').Value = '90.318.94'
$ws.Range('E16').Value = '  -0.84%  '
$ws.Range('D17').Value = '3.728.47'
$ws.Range('E17').Value = '  +1.23%  '
$ws.Range('D18').Value = '3.120.05'
$ws.Range('E18').Value = '  -0.22%  '
$ws.Range('E19').Value = '  -0.07%  '
$ws.Range('D20').Value = '0.0000217'
$ws.Range('E20').Value = '  -1.20%  '
$ws.Range('D21').Value = '14.52'
$ws.Range('E21').Value = '  +3.09%  '
$ws.Range('D22').Value = '449.15'
$ws.Range('E22').Value = '  +3.32%  '
$ws.Range('D23').Value = '5.65'
$ws.Range('E23').Value = '  +9.10%  '
$ws.Range('D24').Value = '9.08'
$ws.Range('E24').Value = '  +2.90%  '
$ws.Range('D25').Value = '6.03'
$ws.Range('E25').Value = '  -2.81%  '
$ws.Range('D26').Value = '90.81'
$ws.Range('E26').Value = '  +5.24%  '
$ws.Range('D27').Value = '12.44'
$ws.Range('E27').Value = '  +1.99%  '
$ws.Range('D28').Value = '3.284.89'
$ws.Range('E28').Value = '  +0.34%  '
$ws.Range('E29').Value = '  +0.14%  '
$ws.Range('D30').Value = '9.72'
$ws.Range('E30').Value = '  +6.67%  '
$ws.Range('E31').Value = '  -4.16%  '
$ws.Range('D32').Value = '27.33'
$ws.Range('E32').Value = '  +15.84%  '
$ws.Range('E33').Value = '  +31.38%  '
$ws.Range('D34').Value = '3.85'
$ws.Range('E34').Value = '  +2.50%  '
$ws.Range('D35').Value = '519.06'
$ws.Range('E35').Value = '  -1.33%  '
$ws.Range('D36').Value = '0.151'
$ws.Range('E36').Value = '  +3.01%  '
$ws.Range('B37').Value = 'RenderToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D37').Value = '7.17'
$ws.Range('E37').Value = '  +0.31%  '
$ws.Range('B38').Value = 'PancakeSwap'
$ws.Range('C38').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D38').Value = '1.95'
$ws.Range('E38').Value = '  +4.70%  '
$ws.Range('E39').Value = '  +0.92%  '
$ws.Range('D40').Value = '0.807'
$ws.Range('E40').Value = '  -10.08%  '
$ws.Range('E41').Value = '  +6.61%  '
$ws.Range('B42').Value = 'Hedera'
$ws.Range('C42').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D42').Value = '0.0865'
$ws.Range('E42').Value = '  -3.60%  '
$ws.Range('B43').Value = 'WhiteBITCoin'
$ws.Range('C43').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D43').Value = '22.18'
$ws.Range('E43').Value = '  -0.36%  '
$ws.Range('E44').Value = '  -0.01%  '
$ws.Range('D45').Value = '3.37'
$ws.Range('E45').Value = '  +44.50%  '
$ws.Range('D46').Value = '1.94'
$ws.Range('E46').Value = '  +1.37%  '
$ws.Range('D47').Value = '0.714'
$ws.Range('E47').Value = '  +14.39%  '
$ws.Range('D48').Value = '151.26'
$ws.Range('E48').Value = '  +2.13%  '
$ws.Range('D49').Value = '46.14'
$ws.Range('E49').Value = '  +5.18%  '
$ws.Range('D50').Value = '4.61'
$ws.Range('E50').Value = '  +8.76%  '
$ws.Range('D51').Value = '1.37'
$ws.Range('E51').Value = '  +4.41%  '
